$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "peter"
$ws.Range("B4").Value = "admin123"
$ws.Range("C4").Value = "Invalid credentials"

$ws.Range("A5").Value = "john"
$ws.Range("B5").Value = "admin123"
$ws.Range("C5").Value = "Invalid credentials"

$excel.Goto($ws.Range("A5:C5"))
